{"js": "// Sprint Retrospective edit: add forgotten team member \"Peter\" in two places,\n// and tidy up a couple of sentences/paragraphs around those mentions.\n\n// 1) \"Marius, Callum, Daniel and Ruth\" -> \"Marius, Callum, Peter, Daniel and Ruth\"\nconst teamListResults = context.document.body.search(\"Marius, Callum, Daniel and Ruth\", { matchCase: true });\nteamListResults.load(\"text\");\nawait context.sync();\nif (teamListResults.items.length > 0) {\n  teamListResults.items[0].insertText(\"Marius, Callum, Peter, Daniel and Ruth\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the trailing space before the end of \"...comply with the clients wishes. \"\nconst wishesResults = context.document.body.search(\"clients wishes. \", { matchCase: true });\nwishesResults.load(\"text\");\nawait context.sync();\nif (wishesResults.items.length > 0) {\n  wishesResults.items[0].insertText(\"clients wishes.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Delete the now-stray empty \"List Paragraph\"-styled paragraph that follows\n//    that sentence (an empty spacer paragraph removed as part of the edit).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"List Paragraph\" && p.text === \"\") {\n    p.delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 4) \"...focus on the database and documentation. \" ->\n//    \"...focus on the database. Ruth and Peter worked on the documentation. \"\nconst ruthDanielResults = context.document.body.search(\"Ruth and Daniel\", { matchCase: true });\nruthDanielResults.load(\"text\");\nawait context.sync();\nif (ruthDanielResults.items.length > 0) {\n  const lastMention = ruthDanielResults.items[ruthDanielResults.items.length - 1];\n  lastMention.insertText(\" focus on the database. Ruth and Peter worked on the \", Word.InsertLocation.after);\n  await context.sync();\n\n  const tailResults = context.document.body.search(\" focus on the database and documentation. \", { matchCase: true });\n  tailResults.load(\"text\");\n  await context.sync();\n  if (tailResults.items.length > 0) {\n    tailResults.items[0].insertText(\"documentation. \", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Sprint Retrospective edit: add forgotten team member \"Peter\" in two places,\n# and tidy up a couple of sentences/paragraphs around those mentions.\n\n$d = $word.ActiveDocument\n\n# 1) \"Marius, Callum, Daniel and Ruth\" -> \"Marius, Callum, Peter, Daniel and Ruth\"\n$find1 = $d.Content.Find\n$find1.Text = \"Marius, Callum, Daniel and Ruth\"\n$find1.Replacement.Text = \"Marius, Callum, Peter, Daniel and Ruth\"\n$find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Remove the trailing space before the end of \"...comply with the clients wishes. \"\n$find2 = $d.Content.Find\n$find2.Text = \"clients wishes. \"\n$find2.Replacement.Text = \"clients wishes.\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 3) Delete the now-stray empty \"List Paragraph\"-styled paragraph that follows\n#    that sentence (an empty spacer paragraph removed as part of the edit).\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Style.NameLocal -eq \"List Paragraph\" -and $p.Range.Text.Trim() -eq \"\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 4) \"...focus on the database and documentation. \" ->\n#    \"...focus on the database. Ruth and Peter worked on the documentation. \"\n# Insert the new text right after the last \"Ruth and Daniel\" mention (keeps the\n# _GoBack bookmark sitting just before \"documentation.\" like in the target doc).\n$find3 = $d.Content.Find\n$find3.Text = \"the programming side whilst Ruth and Daniel\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $ctxRange = $find3.Parent\n    $ctxRange.Collapse(0)\n    $ctxRange.InsertAfter(\" focus on the database. Ruth and Peter worked on the \")\n}\n\n$find4 = $d.Content.Find\n$find4.Text = \" focus on the database and documentation. \"\n$find4.Replacement.Text = \"documentation. \"\n$find4.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
